$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 61 (pushes the existing "27/04/2020".."30/04/2020"
# rows down by one, from rows 61-64 to rows 62-65), then fill in the new
# "01/05/2020" data point.
$ws.Rows.Item(61).Insert()

# Write the new date as literal text (not an auto-converted date serial):
# temporarily mark the cell as Text, assign the value, then restore the
# cell's style to Normal so no stray number-format style is left behind.
$ws.Range("A61").NumberFormatLocal = "@"
$ws.Range("A61").Value = "01/05/2020"
$ws.Range("A61").Style = "Normal"
$ws.Range("B61").Value = 2

# The revised case count for 30/04/2020 (now shifted down to row 65).
$ws.Range("B65").Value = 561
